$d = $word.ActiveDocument

$addition = " We will submit a protocol to the University of Michigan IRB describing our evaluation plan to confirm this exemption."

# Find the paragraph that ends "... observation of public behavior". -
# i.e. the exemption-category paragraph - without relying on literal
# curly-quote characters in this script file.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*exempt under Category 2*") {
        $target = $p
        break
    }
}
if ($target -eq $null) {
    throw "Could not locate the exemption-category paragraph"
}

# The paragraph's very last run is a lone "." immediately before the
# paragraph mark. Collapse a range right after that "." (but before the
# paragraph mark) without selecting/deleting anything, so the preceding
# runs (including the closing curly-quote run) are left completely
# untouched.
$pEnd = $target.Range.End
$insertionPoint = $d.Range($pEnd - 1, $pEnd - 1)

# Insert the new sentence as a brand-new run right after the ".". This
# keeps every earlier run in the paragraph byte-for-byte as it was.
$insertionPoint.InsertAfter($addition + " ")
$newEnd = $insertionPoint.End

# Now rewrite the span covering the old "." run plus the run we just
# inserted (but starting exactly at the "." - never touching the
# curly-quote run before it) with the final combined text. Because this
# text differs from what's currently there, the engine collapses the
# span back down into a single run, giving the same run structure as
# the target edit: the "." run's text becomes
# ". We will submit ... exemption." while the quote run stays separate.
$combined = $d.Range($pEnd - 1, $newEnd)
$combined.Text = $addition
